$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.625.51'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '3.107.77'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'235.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.22%  '

$ws.Range("D6").Value = "'613.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.67%  '

$ws.Range("E7").Value = '  -2.48%  '

$ws.Range("D8").Value = "'0.389"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").Value = '3.104.66'
$ws.Range("E10").Value = '  -0.39%  '

$ws.Range("E11").Value = '  +4.21%  '

$ws.Range("E12").Value = '  -3.74%  '

$ws.Range("E13").Value = '  -3.17%  '

$ws.Range("D14").Value = '92.341.78'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("D15").Value = "'33.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.46%  '

$ws.Range("E16").Value = '  -3.33%  '

$ws.Range("D17").Value = '3.683.23'
$ws.Range("E17").Value = '  -0.66%  '

$ws.Range("D18").Value = '3.105.91'
$ws.Range("E18").Value = '  -2.59%  '

$ws.Range("D19").Value = "'3.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("D20").Value = "'14.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.06%  '

$ws.Range("D21").Value = "'5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.56%  '

$ws.Range("E22").Value = '  +1.12%  '

$ws.Range("E23").Value = '  -4.10%  '

$ws.Range("D24").Value = "'9.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").Value = "'8.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.38%  '

$ws.Range("E26").Value = '  -6.52%  '

$ws.Range("D27").Value = "'85.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.75%  '

$ws.Range("D28").Value = "'11.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.60%  '

$ws.Range("D29").Value = '3.273.60'
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("E31").Value = '  +8.43%  '

$ws.Range("D32").Value = "'0.234"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.93%  '

$ws.Range("E33").Value = '  -12.38%  '

$ws.Range("D34").Value = "'9.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.54%  '

$ws.Range("E35").Value = '  -28.30%  '

$ws.Range("D36").Value = "'8.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.39%  '

$ws.Range("D37").Value = "'0.166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.34%  '

$ws.Range("D38").Value = "'25.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.06%  '

$ws.Range("D39").Value = "'4.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.43%  '

$ws.Range("E40").Value = '  -15.06%  '

$ws.Range("E41").Value = '  +7.75%  '

$ws.Range("E42").Value = '  -2.62%  '

$ws.Range("D43").Value = "'463.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.48%  '

$ws.Range("D44").Value = "'0.428"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.35%  '

$ws.Range("D45").Value = "'3.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.31%  '

$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("D47").Value = "'159.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '

$ws.Range("D48").Value = "'0.680"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.81%  '

$ws.Range("E49").Value = '  -4.85%  '

$ws.Range("E50").Value = '  -2.68%  '

$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = "'43.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
